# Rename column headers from "Title Case" labels to camelCase keys, and
# let the column widths follow the now-shorter header text (mirrors
# Excel's behaviour of the columns having been auto-fit to the new text).

$wb = $excel.ActiveWorkbook

# Offset between the `ColumnWidth` property (measured in characters) and
# the raw `width` value stored in the worksheet XML for the Calibri 11
# default font used by this workbook.
$widthOffset = 5 / 6

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1").Value = "usage"
    $ws.Range("B1").Value = "viewportWidth"
    $ws.Range("C1").Value = "pixelRatio"
    $ws.Range("D1").Value = "imgWidth"
    $ws.Range("E1").Value = "imgVW"
    $ws.Range("F1").Value = "idealIntrinsicWidth"
    $ws.Range("G1").Value = "chosenIntrinsicWidth"
    $ws.Range("H1").Value = "renderedFidelity"
    $ws.Range("I1").Value = "renderedToIdealFidelityRatio"
    $ws.Range("J1").Value = "evaluation"
    $ws.Range("K1").Value = "waste"

    $ws.Range("B1").ColumnWidth = 14 - $widthOffset
    $ws.Range("C1").ColumnWidth = 11 - $widthOffset
    $ws.Range("D1").ColumnWidth = 9 - $widthOffset
    $ws.Range("E1").ColumnWidth = 6 - $widthOffset
    $ws.Range("F1").ColumnWidth = 20 - $widthOffset
    $ws.Range("G1").ColumnWidth = 21 - $widthOffset
    $ws.Range("H1").ColumnWidth = 17 - $widthOffset
    $ws.Range("I1").ColumnWidth = 29 - $widthOffset
}
